$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (ownTeam, oppTeam), shifting old
# batsman/totalRuns/totalBalls/total4s/total6s/sr columns from D..I to F..K
$ws.Range("D1:E1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Row 2 - fill in the new ownTeam / oppTeam values for the existing record
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"

# New row 3
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 01 2020"
$ws.Range("C3").Value = "Mumbai won by 48 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "James Neesham "
$ws.Range("G3").Value = "'7"
$ws.Range("H3").Value = "'7"
$ws.Range("I3").Value = "'0"
$ws.Range("J3").Value = "'0"
$ws.Range("K3").Value = "'100.00"

# New row 4
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " October 20 2020"
$ws.Range("C4").Value = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Delhi Capitals"
$ws.Range("F4").Value = "James Neesham "
$ws.Range("G4").Value = "'10"
$ws.Range("H4").Value = "'8"
$ws.Range("I4").Value = "'0"
$ws.Range("J4").Value = "'1"
$ws.Range("K4").Value = "'125.00"
